$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# The report for file "f02b04e1-c7d8-443c-95f0-c0079bb6c5f1.md" (row 3 on every
# sheet) has now been handed back and is in sync with en-US, so update the
# status columns and the "Latest Handback DateTime" columns accordingly.

$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-03-10 20:53:14"

$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-03-10 20:53:26"
